$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the GREY header (was "GREY`n(Yes/No)") to just "GREY"
$ws.Range("AG1").Value = "GREY"

# Add the Grey value for the data row
$ws.Range("AG2").Value = "Yes"
